$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1718
$ws1.Range("F6").Value = 3294
$ws1.Range("F7").Value = 949
$ws1.Range("F8").Value = 2122
$ws1.Range("F9").Value = 2051
$ws1.Range("F10").Value = 1065
$ws1.Range("F13").Value = 1639
$ws1.Range("F16").Value = 21
$ws1.Range("F18").Value = 133
$ws1.Range("F19").Value = 1505
$ws1.Range("F20").Value = 569
$ws1.Range("F21").Value = 668
$ws1.Range("F22").Value = 555
$ws1.Range("F23").Value = 11952
$ws1.Range("F24").Value = 11967
$ws1.Range("F25").Value = 879
$ws1.Range("F26").Value = 675
$ws1.Range("F28").Value = 1875
$ws1.Range("F30").Value = 499

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 7

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1718
$ws4.Range("F8").Value = 3294
$ws4.Range("F9").Value = 949
$ws4.Range("F10").Value = 2122
$ws4.Range("F11").Value = 2051
$ws4.Range("F12").Value = 1065
$ws4.Range("F15").Value = 1639
$ws4.Range("F18").Value = 21
$ws4.Range("F22").Value = 133
$ws4.Range("F23").Value = 1505
$ws4.Range("F24").Value = 569
$ws4.Range("F25").Value = 668
$ws4.Range("F26").Value = 555
$ws4.Range("F27").Value = 11952
$ws4.Range("F28").Value = 11967
$ws4.Range("F29").Value = 879
$ws4.Range("F30").Value = 675
$ws4.Range("F32").Value = 1875
$ws4.Range("F36").Value = 499
$ws4.Range("F37").Value = 7
